$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.144070506095886
$ws.Range("B1").Value = 2.29717493057251
$ws.Range("D1").Value = 1.470521926879883
$ws.Range("E1").Value = 0.9556595087051392
